# Applies the Ifrit_Profits crafting-leve price/profit refresh captured in the
# scheduled-runner commit: updates cached currentAveragePrice* / LevePrice* /
# LeveProfit* figures (cols H:N) for the affected leve rows across the ALC,
# ARM, BSM, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

# Row 53 on ALC (Leve Item ID 5479)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 950.4286
$ws.Range("I53").Value = 291.5
$ws.Range("J53").Value = 1829
$ws.Range("K53").Value = 291.5
$ws.Range("L53").Value = 1829
$ws.Range("M53").Value = 345.5
$ws.Range("N53").Value = -3103

# Row 106 on ALC (Leve Item ID 19903)
$ws.Range("H106").Value = 2737
$ws.Range("I106").Value = 2776.5
$ws.Range("J106").Value = 2500
$ws.Range("K106").Value = 2776.5
$ws.Range("L106").Value = 2500
$ws.Range("M106").Value = -2145.5
$ws.Range("N106").Value = -3762

# Row 132 on ALC (Leve Item ID 44049)
$ws.Range("H132").Value = 373889.2
$ws.Range("I132").Value = 388192.62
$ws.Range("K132").Value = 1164577.86
$ws.Range("M132").Value = -1162047.86

# Row 61 on ARM (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2340.8
$ws.Range("I61").Value = 1629.0769
$ws.Range("J61").Value = 3111.8333
$ws.Range("K61").Value = 1629.0769
$ws.Range("L61").Value = 3111.8333
$ws.Range("M61").Value = -1417.0769
$ws.Range("N61").Value = -3535.8333

# Row 74 on ARM (Leve Item ID 44000)
$ws.Range("H74").Value = 3157.6667
$ws.Range("I74").Value = 690.4783
$ws.Range("J74").Value = 4988.161
$ws.Range("K74").Value = 690.4783
$ws.Range("L74").Value = 4988.161
$ws.Range("M74").Value = 183.5217
$ws.Range("N74").Value = -6736.161

# Row 77 on ARM (Leve Item ID 44000)
$ws.Range("H77").Value = 3157.6667
$ws.Range("I77").Value = 690.4783
$ws.Range("J77").Value = 4988.161
$ws.Range("K77").Value = 3452.3915
$ws.Range("L77").Value = 24940.805
$ws.Range("M77").Value = 915.6085000000003
$ws.Range("N77").Value = -33676.805

# Row 97 on ARM (Leve Item ID 19941)
$ws.Range("H97").Value = 454.4
$ws.Range("I97").Value = 456.35715
$ws.Range("J97").Value = 449.83334
$ws.Range("K97").Value = 456.35715
$ws.Range("L97").Value = 449.83334
$ws.Range("M97").Value = 39.64285000000001
$ws.Range("N97").Value = -1441.83334

# Row 102 on ARM (Leve Item ID 19945)
$ws.Range("H102").Value = 1488.6471
$ws.Range("I102").Value = 1522.6666
$ws.Range("K102").Value = 1522.6666
$ws.Range("M102").Value = 99.33339999999998

# Row 110 on ARM (Leve Item ID 27708)
$ws.Range("H110").Value = 1007.7143
$ws.Range("I110").Value = 903.4583
$ws.Range("K110").Value = 903.4583
$ws.Range("M110").Value = 1141.5417

# Row 122 on ARM (Leve Item ID 36168)
$ws.Range("H122").Value = 1670.9333
$ws.Range("I122").Value = 1729.9166
$ws.Range("J122").Value = 1435
$ws.Range("K122").Value = 5189.7498
$ws.Range("L122").Value = 4305
$ws.Range("M122").Value = -2739.7498
$ws.Range("N122").Value = -9205

# Row 136 on ARM (Leve Item ID 43999)
$ws.Range("H136").Value = 2340.8
$ws.Range("I136").Value = 1629.0769
$ws.Range("J136").Value = 3111.8333
$ws.Range("K136").Value = 4887.2307
$ws.Range("L136").Value = 9335.499899999999
$ws.Range("M136").Value = -2337.2307
$ws.Range("N136").Value = -14435.4999

# Row 99 on BSM (Leve Item ID 19943)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 714.25
$ws.Range("I99").Value = 664.93335
$ws.Range("K99").Value = 664.93335
$ws.Range("M99").Value = 833.06665

# Row 105 on BSM (Leve Item ID 19947)
$ws.Range("H105").Value = 1892.8572
$ws.Range("I105").Value = 1892.8572
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1892.8572
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -145.8571999999999
$ws.Range("N105").ClearContents()

# Row 56 on CUL (Leve Item ID 10146)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 1354.6154
$ws.Range("I56").Value = 1354.6154
$ws.Range("K56").Value = 1354.6154
$ws.Range("M56").Value = -824.6153999999999

# Row 107 on CUL (Leve Item ID 27838)
$ws.Range("H107").Value = 77328.66
$ws.Range("I107").Value = 125208.5
$ws.Range("J107").Value = 56048.723
$ws.Range("K107").Value = 375625.5
$ws.Range("L107").Value = 168146.169
$ws.Range("M107").Value = -373705.5
$ws.Range("N107").Value = -171986.169

# Row 122 on CUL (Leve Item ID 36078)
$ws.Range("H122").Value = 15152895
$ws.Range("I122").Value = 20833838
$ws.Range("J122").Value = 3717.3333
$ws.Range("K122").Value = 187504542
$ws.Range("L122").Value = 33455.9997
$ws.Range("M122").Value = -187502092
$ws.Range("N122").Value = -38355.9997

# Row 131 on CUL (Leve Item ID 36060)
$ws.Range("H131").Value = 2177127.8
$ws.Range("J131").Value = 2780124.8
$ws.Range("L131").Value = 8340374.399999999
$ws.Range("N131").Value = -8350454.399999999

# Row 132 on CUL (Leve Item ID 43972)
$ws.Range("H132").Value = 66667748
$ws.Range("I132").Value = 111111816
$ws.Range("K132").Value = 1000006344
$ws.Range("M132").Value = -1000003814

# Row 58 on GSM (Leve Item ID 4363)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 14950
$ws.Range("J58").Value = 14950
$ws.Range("L58").Value = 14950
$ws.Range("N58").Value = -15504

# Row 80 on GSM (Leve Item ID 12521)
$ws.Range("H80").Value = 134998.44
$ws.Range("I80").Value = 2516
$ws.Range("J80").Value = 300601.5
$ws.Range("K80").Value = 2516
$ws.Range("L80").Value = 300601.5
$ws.Range("M80").Value = -1518
$ws.Range("N80").Value = -302597.5

# Row 83 on GSM (Leve Item ID 12521)
$ws.Range("H83").Value = 134998.44
$ws.Range("I83").Value = 2516
$ws.Range("J83").Value = 300601.5
$ws.Range("K83").Value = 12580
$ws.Range("L83").Value = 1503007.5
$ws.Range("M83").Value = -7588
$ws.Range("N83").Value = -1512991.5

# Row 97 on GSM (Leve Item ID 19940)
$ws.Range("H97").Value = 1153.5483
$ws.Range("I97").Value = 1228.2609
$ws.Range("J97").Value = 938.75
$ws.Range("K97").Value = 1228.2609
$ws.Range("L97").Value = 938.75
$ws.Range("M97").Value = -732.2609
$ws.Range("N97").Value = -1930.75

# Row 107 on GSM (Leve Item ID 27802)
$ws.Range("H107").Value = 560.4545000000001
$ws.Range("I107").Value = 345.81818
$ws.Range("J107").Value = 775.0909
$ws.Range("K107").Value = 345.81818
$ws.Range("L107").Value = 775.0909
$ws.Range("M107").Value = 1574.18182
$ws.Range("N107").Value = -4615.0909

# Row 113 on GSM (Leve Item ID 27710)
$ws.Range("H113").Value = 2330.4443
$ws.Range("I113").Value = 5864
$ws.Range("J113").Value = 1320.8572
$ws.Range("K113").Value = 5864
$ws.Range("L113").Value = 1320.8572
$ws.Range("M113").Value = -3694
$ws.Range("N113").Value = -5660.8572

# Row 132 on GSM (Leve Item ID 44008)
$ws.Range("H132").Value = 3098.25
$ws.Range("I132").Value = 2546.6428
$ws.Range("K132").Value = 7639.928400000001
$ws.Range("M132").Value = -5109.928400000001

# Row 135 on GSM (Leve Item ID 42006)
$ws.Range("H135").Value = 33066.5
$ws.Range("J135").Value = 33066.5
$ws.Range("L135").Value = 33066.5
$ws.Range("N135").Value = -43206.5

# Row 55 on LTW (Leve Item ID 5284)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 216.94118
$ws.Range("J55").Value = 187.5
$ws.Range("L55").Value = 187.5
$ws.Range("N55").Value = -533.5

# Row 136 on LTW (Leve Item ID 44060)
$ws.Range("H136").Value = 1523.4706
$ws.Range("I136").Value = 1030
$ws.Range("J136").Value = 2228.4285
$ws.Range("K136").Value = 3090
$ws.Range("L136").Value = 6685.2855
$ws.Range("M136").Value = -540
$ws.Range("N136").Value = -11785.2855

# Row 136 on WVR (Leve Item ID 44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 9368.458000000001
$ws.Range("I136").Value = 21019.8
$ws.Range("J136").Value = 1046.0714
$ws.Range("K136").Value = 63059.39999999999
$ws.Range("L136").Value = 3138.2142
$ws.Range("M136").Value = -60509.39999999999
$ws.Range("N136").Value = -8238.2142
